$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds header "K" (previously based on "Strike#"); regenerate the
# computed K values for each data row (2-12).
$kValues = @{
    2  = 1
    3  = 0
    4  = 4
    5  = 0
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 1
    12 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
